# Generate Report for Handoff
#
# Updates the localization-status workbook to reflect a new handoff run:
#   - "Latest Handoff Datetime" timestamps are refreshed for the files that
#     were just handed off (zh-cn and de-de target files), and the matching
#     "Latest HO Xliff Generate Date" column on the Overview sheet is kept
#     in sync.
#   - The "Priority" column for those same rows is marked "ht" (handed off)
#     for both the zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows 7, 8, 9, 10, 12 and 13 correspond to the files included in this
# handoff batch (row 11 was handed off separately and is left untouched).
$handoffRows = 7,8,9,10,12,13

$zhCnHandoffTime = "2016-09-03 22:23:15"
$deDeHandoffTime = "2016-09-03 22:23:20"

foreach ($r in $handoffRows) {
    # zh-cn worksheet: refresh handoff timestamp and mark priority as handed off
    $wsZhCn.Range("H$r").Value = $zhCnHandoffTime
    $wsZhCn.Range("E$r").Value = "ht"

    # de-de worksheet: refresh handoff timestamp and mark priority as handed off
    $wsDeDe.Range("H$r").Value = $deDeHandoffTime
    $wsDeDe.Range("E$r").Value = "ht"

    # Overview worksheet: keep the "Latest HO Xliff Generate Date" column
    # in sync with the de-de handoff timestamp used above.
    $wsOverview.Range("G$r").Value = $deDeHandoffTime
}
